$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns for the sales-order filter test case (W1:AA1) ---
# Copy the existing bold-Arial header style (used by L1:V1, style index 7)
# onto the new header cells before setting their text.
$ws.Range("S1").Copy() | Out-Null
$ws.Range("W1:AA1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("W1").Value = "ResellerName"
$ws.Range("X1").Value = "EndUserName"
$ws.Range("Y1").Value = "CreatedOn"
$ws.Range("Z1").Value = "FilterOrderType"
$ws.Range("AA1").Value = "FilterOrderStatus"

# --- sales_orders row (row 3): update / add the new filter test data ---
$ws.Range("C3").Value = "20-VN2W9-11"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 400
$ws.Range("W3").Value = "INGRAM MICRO CAP TEST ACCOUNT"
$ws.Range("X3").Value = "Everest EndUser"
$ws.Range("Y3").Value = "Yesterday"
$ws.Range("Z3").Value = "Stock"
$ws.Range("AA3").Value = "Order Hold"

# --- Column V got a touch narrower once the new columns landed ---
$ws.Columns("V").ColumnWidth = 26

# --- Restore the view: scroll back to A1 and leave the selection on F9 ---
$ws.Range("F9").Select() | Out-Null
